# Updated cryptos list on Mon Jan 22 09:44:30 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns with the latest scrape,
# and EnergySwap/VeChain (rows 44-45) swap ranking positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'40.911.92"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "'2.407.28"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'314.47"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").Value = "'88.10"
$ws.Range("E6").Value = "  -5.32%  "
$ws.Range("D7").Value = "'0.535"
$ws.Range("E7").Value = "  -3.15%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "'0.493"
$ws.Range("E9").Value = "  -4.57%  "
$ws.Range("D10").Value = "'0.0829"
$ws.Range("E10").Value = "  -3.27%  "
$ws.Range("D11").Value = "'31.23"
$ws.Range("E11").Value = "  -5.70%  "
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").Value = "'2.780.86"
$ws.Range("E13").Value = "  -2.62%  "
$ws.Range("D14").Value = "'6.73"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").Value = "'15.38"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").Value = "'2.410.00"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("D17").Value = "'0.766"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").Value = "'40.804.82"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").Value = "0.0₃0917"
$ws.Range("E19").Value = "  -3.66%  "
$ws.Range("D20").Value = "'6.21"
$ws.Range("E20").Value = "  -4.20%  "
$ws.Range("D21").Value = "'70.46"
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").Value = "'10.81"
$ws.Range("E22").Value = "  -4.50%  "
$ws.Range("D23").Value = "'237.54"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "'2.65"
$ws.Range("E24").Value = "  -3.36%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  -5.03%  "
$ws.Range("D27").Value = "'23.85"
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("E28").Value = "  -2.67%  "
$ws.Range("D29").Value = "'9.46"
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("D30").Value = "'33.99"
$ws.Range("E30").Value = "  -6.01%  "
$ws.Range("D31").Value = "'156.81"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  -4.98%  "
$ws.Range("E34").Value = "  -4.25%  "
$ws.Range("E35").Value = "  -4.80%  "
$ws.Range("D36").Value = "'2.85"
$ws.Range("E36").Value = "  -3.16%  "
$ws.Range("E37").Value = "  -1.69%  "
$ws.Range("D38").Value = "'16.09"
$ws.Range("E38").Value = "  -7.19%  "
$ws.Range("E39").Value = "  -7.33%  "
$ws.Range("D40").Value = "'0.0993"
$ws.Range("E40").Value = "  -4.38%  "
$ws.Range("D41").Value = "'3.84"
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("E42").Value = "  -6.92%  "
$ws.Range("D43").Value = "'1.983.01"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0272"
$ws.Range("E44").Value = "  -4.73%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'18.01"
$ws.Range("E45").Value = "  -5.42%  "
$ws.Range("D46").Value = "'2.83"
$ws.Range("E46").Value = "  -5.04%  "
$ws.Range("D47").Value = "'9.32"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "'2.649.54"
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("D49").Value = "'73.77"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").Value = "'93.58"
$ws.Range("E50").Value = "  -3.86%  "
$ws.Range("D51").Value = "'50.98"
$ws.Range("E51").Value = "  -2.47%  "
